$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "[-, -, 'MEC-1NA-Elet. Dig. Bas.', -]"

$ws.Range("F18").Value = $newValue
$ws.Range("F19").Value = $newValue
$ws.Range("F20").Value = $newValue
$ws.Range("F21").Value = $newValue
